$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "start"/"end" header labels between columns B and C
$ws.Range("B1").Value = "start"
$ws.Range("C1").Value = "end"
